$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A17").Value = "a"
$ws.Range("B17").Value = 0

$ws.Range("A18").Value = "a"
$ws.Range("B18").Value = 1042
